$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 ("「波」moug" entry). All following rows shift up by one.
$ws.Rows.Item(9).Delete()
